$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

# New Neurolucida results DEG run-log entries
$ws.Cells.Item(107,1).Value = "Neurolucida results"
$ws.Cells.Item(107,2).Value = "2022-07-04 11-56-43"
$ws.Cells.Item(107,3).Value = "DEG"
$ws.Cells.Item(107,4).Value = "A"

$ws.Cells.Item(108,1).Value = "Neurolucida results"
$ws.Cells.Item(108,2).Value = "2022-07-04 11-58-43"
$ws.Cells.Item(108,3).Value = "DEG"
$ws.Cells.Item(108,4).Value = "C"

$ws.Cells.Item(109,1).Value = "Neurolucida results"
$ws.Cells.Item(109,2).Value = "2022-07-04 11-59-22"
$ws.Cells.Item(109,3).Value = "DEG"
$ws.Cells.Item(109,4).Value = "N"

$null = $ws.Range("B111").Select()
